$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 2
$ws.Range("G8").Value = 2
$ws.Range("G9").Value = 3
$ws.Range("G10").Value = 5
